$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank "G" (Saturday) attendance column with the
# newly submitted attendance marks for 3 October.
$ws.Range("G10").Value = 28

$ws.Range("G12").Value = 3

$ws.Range("G14").Value = 3
$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 3
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 3

# Reset the view: scroll back to the top-left and move the active selection
# from the merged "Total" row down to I18.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I18").Select()
